$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Title
Replace-Text "Circuitry of Human Creativity: Unlocking Inspiration" "Exploring the Marvels of the Human Body: An Immersive Journey into Biology's Wonders"

# Author name (merges "Dr" + "." + " Beatrice Carter" into a single run)
Replace-Text "Dr. Beatrice Carter" "Mark Lawson"

# Email address parts
Replace-Text "beatrice" "mark"
Replace-Text "carter@interlink" "lawson@educationalhaven"

# Body paragraph 1 (first group of three sentences)
Replace-Text "The human brain, an intricate tapestry of neurons and synapses, harbors an enigmatic phenomenon known as creativity: the ability to generate novel ideas, concepts, and solutions" "Biology, a realm of astounding complexities and awe-inspiring discoveries, unveils the intricacies of life on our planet"

Replace-Text " While often shrouded in mystery, creativity presents an alluring frontier for researchers seeking to uncover the mechanisms that underlie this extraordinary cognitive process" " It weaves together the tapestry of living organisms, unlocking their secrets and revealing the marvels of existence"

Replace-Text " Recent advancements in neuroscience, psychology, and artificial intelligence have illuminated the intricate circuitry of human creativity, revealing a complex interplay between brain regions, cognitive processes, and life experiences" " From the smallest microorganisms to the grandest ecosystems, biology captivates our imagination and challenges us to explore the depths of life's mysteries"

# Body paragraph 1 (second group)
Replace-Text "The birth of a novel idea, whether a scientific breakthrough, artistic masterpiece, or innovative solution, can often feel as unexpected as a bolt of lightning" "Biology's profound impact extends beyond mere scientific knowledge; it enriches our understanding of ourselves, our place in the universe, and our interconnectedness with all living beings"

# This merges " However, beneath...diverse sources" + "." + " The prefrontal cortex...new concepts" into one run
Replace-Text " However, beneath this seemingly random flash of inspiration, complex cognitive machinery churns, analyzing, synthesizing, and recombining information from diverse sources. The prefrontal cortex, a region associated with higher-order cognitive processes, serves as a central hub for creativity, facilitating the integration of diverse stimuli and the generation of new concepts" " It unravels the enigma of human development, shedding light on the intricacies of our bodies, the wonders of our minds, and the profound symphony of our genetic heritage"

# Body paragraph 1 (third group)
Replace-Text "Neurotransmitters, the chemical messengers of the brain, play a crucial role as facilitators and modulators of creativity" "Our exploration of biology unveils the harmony of life's processes, akin to a captivating quantum dance"

Replace-Text " Dopamine, known for its salience in reward pathways, prompts the exploration of novel ideas and encourages risk-taking" " The interdependence of organisms, the delicate balance of ecosystems, and the remarkable resilience of life amidst adversity evoke a sense of awe and wonder"

# This merges " Opioid systems...of inspiration" + "." + " Our experiences, learning...unfolds" into one run
Replace-Text " Opioid systems provide internal validation for creative thinking, reinforcing and rewarding moments of inspiration. Our experiences, learning, and environment also shape the creative landscape of our minds, providing both the raw material and the context in which creativity unfolds" " Each revelation, each puzzle solved, propels us further into the depths of this magnificent science"

# Summary heading paragraph stays the same ("Summary")

# Summary body
Replace-Text "The circuitry of human creativity is a intricate tapestry of brain regions, cognitive processes, and life experiences" "Biology unveils the astounding complexities and awe-inspiring discoveries hidden within the realm of life on Earth"

Replace-Text " The prefrontal cortex serves as a central hub for creativity, facilitating the integration of diverse stimuli and the generation of new concepts" " It invites us to explore the intricacies of living organisms, unraveling the marvels of existence"

# This merges " Neurotransmitters " + "act as chemical messengers...thinking" (removing the lastRenderedPageBreak run) into one run
Replace-Text " Neurotransmitters act as chemical messengers, influencing our ability to explore novel ideas and providing internal validation for creative thinking" " Biology's impact extends beyond scientific knowledge, enriching our understanding of ourselves, our place in the universe, and our interconnectedness with all living beings"

# This merges " Our experiences, ranging...creativity" + "." + " Thus, understanding...innovators" into one run
Replace-Text " Our experiences, ranging from exposure to the arts to cultural influences and personal struggles, shape the fabric of our creativity. Thus, understanding the circuitry of creativity can not only unlock the potential of human ingenuity but also pave the way for fostering creativity in diverse fields and nurturing the minds of future innovators" " Through its revelations, biology evokes a sense of awe and wonder, propelling us further into the depths of this magnificent science"

# Add a new empty paragraph at the end of the document body (before sectPr)
$end = $d.Content
$end.Collapse(0) | Out-Null
$end.InsertParagraphAfter() | Out-Null
